$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Mark rows 43-46 (column C/D) as "NA" (previously blank), center-aligned like the
#    other Y/N cells in those columns (style index 3 == horizontal center alignment).
foreach ($r in 43..46) {
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = "NA"
    $cCell.HorizontalAlignment = -4108   # xlCenter

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = "NA"
    $dCell.HorizontalAlignment = -4108   # xlCenter
}

# 2. Insert a new row before the current row 232 (WriteTimeSeriesToHydroJSON) to hold the
#    new WriteTimeSeriesToGeoJSON command documentation, keeping alphabetical order.
$ws.Rows.Item(232).Insert()

$ws.Cells.Item(232, 1).Value = "WriteTimeSeriesToGeoJSON"
$ws.Cells.Item(232, 2).Value = "Write time series to GeoJSON file."
$ws.Cells.Item(232, 3).Value = "Y"
$ws.Cells.Item(232, 4).Value = "Y"

# 3. Update the summary COUNTIF formulas (now on row 237 after the insert) to also count
#    "NA" entries in addition to "Y" entries.
$ws.Cells.Item(237, 3).Formula = '=COUNTIF(C2:C236,"=Y")+COUNTIF(C2:C236,"=NA")'
$ws.Cells.Item(237, 4).Formula = '=COUNTIF(D2:D236,"=Y")+COUNTIF(D2:D236,"=NA")'

# 4. Update the active/frozen-pane selection to reflect the new bottom of the sheet.
$ws.Range("B219").Select()
$ws.Range("D238").Select()
